# Katalog guncellendi - Sal 25.11.2025 10:53:13,17
# Adds 6 new "ERKEK NUBUK CEKET" product rows to the bottom of the catalogue.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fiyat = "440 TL"
$kategori = "Ceket"
$aciklama = "S-M-L-XL-2XL Beden seçeneği mevcuttur.Dilerseniz battal beden mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır."
$stok = "Var"

# First new row (89) is typed in natural left-to-right column order.
$ws.Cells.Item(89, 1).Value = "ERKEK NUBUK CEKET TABA"
$ws.Cells.Item(89, 2).Value = $fiyat
$ws.Cells.Item(89, 3).Value = $kategori
$ws.Cells.Item(89, 4).Value = "TABANUBUK.jpg"
$ws.Cells.Item(89, 5).Value = $aciklama
$ws.Cells.Item(89, 6).Value = $stok

# Remaining colour variants (90-94): the image filename (D) was filled in
# before the product name (A) for each row.
$varyantlar = @(
    @{ Gorsel = "GRİNUBUK.jpg";        Ad = "ERKEK NUBUK CEKET GRİ" },
    @{ Gorsel = "HAKİNUBUK.jpg";       Ad = "ERKEK NUBUK CEKET HAKİ" },
    @{ Gorsel = "KAHVERENGİNUBUK.jpg"; Ad = "ERKEK NUBUK CEKET KAHVERENGİ" },
    @{ Gorsel = "LACİVERTNUBUK.jpg";   Ad = "ERKEK NUBUK CEKET LACİVERT" },
    @{ Gorsel = "SİYAHNUBUK.jpg";      Ad = "ERKEK NUBUK CEKET SİYAH" }
)

$row = 90
foreach ($varyant in $varyantlar) {
    $ws.Cells.Item($row, 4).Value = $varyant.Gorsel
    $ws.Cells.Item($row, 1).Value = $varyant.Ad
    $ws.Cells.Item($row, 2).Value = $fiyat
    $ws.Cells.Item($row, 3).Value = $kategori
    $ws.Cells.Item($row, 5).Value = $aciklama
    $ws.Cells.Item($row, 6).Value = $stok
    $row = $row + 1
}

$lastRow = $row - 1
$excel.ActiveWindow.ScrollRow = 73
$ws.Range("E" + $lastRow).Select()
